# B1/B2 PowerPoint update
# 1. Re-style the slide 5 table with the alternate ("Office"-family) table
#    style instead of the custom "Table_0" style.
# 2. Swap the deck's colour theme: the slide master currently carries the
#    "Integral" (Red Violet) theme colours; replace them with the stock
#    "Office Theme" (Office) colours - i.e. apply the Office colour theme
#    to the deck.

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------------
$tableSlide = $p.Slides.Item(5)
for ($i = 1; $i -le $tableSlide.Shapes.Count; $i++) {
    $shp = $tableSlide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{DA613774-7D5C-4B3F-8359-748A3C7ABF42}")
    }
}

# --- 2. Theme colours -------------------------------------------------------
# MsoThemeColorSchemeIndex order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeThemeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$colorScheme = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $colorScheme.Item($i).RGB = $officeThemeColors[$i - 1]
}
